# Apply the commit: "removed ER tags from non-ER templates and non-ER tags"
#
# 1. Rename the "SwateTemplateMetadata" sheet to "isa_template".
# 2. On that sheet, the "Tags" list (row 12), "Tags Term Accession Number"
#    (row 13) and "Tags Term Source REF" (row 14) currently repeat the ER
#    value ("GEO" / "DPBO:1000096" / "DPBO") as the first tag entry. Remove
#    that duplicated first entry and shift the remaining tag values one
#    column to the left, clearing the now-unused trailing column.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("SwateTemplateMetadata")
$ws.Name = "isa_template"

# Rows 12-14 hold the "Tags" list (Tags / Tags Term Accession Number / Tags
# Term Source REF), each currently starting with a duplicate of the single
# ER entry (GEO / DPBO:1000096 / DPBO). Shift C:E one column left into B:D
# (copy carries the formatting/style along, matching how the remaining tag
# entries - and the trailing blank list cells - looked before the removed
# one) and drop the now-superfluous source cells.
$ws.Range("C12:E14").Copy($ws.Range("B12:D14"))

# The copy duplicates C13/C14's old values into themselves; clear those back
# out now that B13/B14 carry them instead.
$ws.Range("C13").ClearContents()
$ws.Range("C14").ClearContents()

# Column E in rows 12-14 is now unused (the lists only span B:D) - delete the
# cells outright (dropping column E from the sheet's used range) rather than
# merely clearing their contents.
$ws.Range("E12:E14").Delete() | Out-Null

# Row 12 no longer needs its custom (wrapped-text) height now that it holds
# one fewer entry - let it size back down to the default row height.
$ws.Rows(12).AutoFit()
